$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = '[{"factor": ["thioredoxin concentration", "serum"], "outcome": ["in-hospital major adverse event"]}, {"factor": ["admission", "glasgow coma scale score"], "outcome": ["in-hospital major adverse event"]}]'
$ws.Range("G5").Value = '[{"factor": ["glasgow coma scale score"], "outcome": ["glasgow outcome scale", "gos"]}, {"factor": ["artificial ventilation"], "outcome": ["mortality rate"]}, {"factor": ["intracranial", "monitoring"], "outcome": ["long-term", "outcome", "good"]}]'
$ws.Range("G6").Value = '[{"factor": ["impact", "calculator", "prognosis"], "outcome": ["elderly", "tbi", "severe", "outcome", "patient"]}, {"factor": ["expect risk", "outcome", "fatal"], "outcome": ["risk", "conservative treatment"]}, {"factor": ["unfavorable", "outcome", "predict risk"], "outcome": ["risk", "outcome", "rate", "actual", "predict", "unfavorable"]}]'
$ws.Range("G7").Value = '[{"factor": ["time post-tbi"], "outcome": ["percentage", "cognitive functionality gain"]}, {"factor": ["admission", "cognitive function"], "outcome": ["percentage", "cognitive functionality gain"]}]'
$ws.Range("G8").Value = '[{"factor": ["ct", "compute tomography"], "outcome": ["patient", "tbi", "death"]}, {"factor": ["marshall and rotterdam scoring system"], "outcome": ["patient", "tbi", "death"]}, {"factor": ["basal cistern absence"], "outcome": ["patient", "tbi", "death"]}, {"factor": ["positive midline shift"], "outcome": ["patient", "tbi", "death"]}, {"factor": ["hemorrhagic mass volume"], "outcome": ["patient", "tbi", "death"]}, {"factor": ["intraventricular", "subarachnoid hemorrhage"], "outcome": ["patient", "tbi", "death"]}]'
$ws.Range("G9").Value = '[{"factor": ["apache ii"], "outcome": ["icu-treated", "tbi", "six-month", "patient", "mortality"]}, {"factor": ["sap ii"], "outcome": ["icu-treated", "tbi", "six-month", "patient", "mortality"]}, {"factor": ["sofa"], "outcome": ["icu-treated", "tbi", "six-month", "patient", "mortality"]}, {"factor": ["age"], "outcome": ["icu-treated", "tbi", "six-month", "patient", "mortality"]}, {"factor": ["glasgow coma scale"], "outcome": ["icu-treated", "tbi", "six-month", "patient", "mortality"]}]'
$ws.Range("G10").Value = '[{"factor": ["v/c ratio"], "outcome": ["gos score", "lcf score", "drs"]}]'
$ws.Range("G11").Value = '[{"factor": ["timp-1 level", "serum"], "outcome": ["patient", "tbi", "mortality", "severe"]}]'
$ws.Range("G12").Value = '[{"factor": ["contusion", "mri", "evidence"], "outcome": ["glasgow outcome scale-extended", "gos-e"]}, {"factor": ["reduce", "roi", "severely", "fa"], "outcome": ["gos-e"]}, {"factor": ["neuropsychiatric history"], "outcome": ["gos-e"]}, {"factor": ["age"], "outcome": ["gos-e"]}, {"factor": ["year of"], "outcome": ["gos-e"]}, {"factor": ["reduce", "roi", "severely", "fa"], "outcome": ["gos-e"]}, {"factor": ["neuropsychiatric history"], "outcome": ["gos-e"]}, {"factor": ["year of"], "outcome": ["gos-e"]}]'
$ws.Range("G13").Value = '[{"factor": ["depressive symptom", "preinjury"], "outcome": ["health-related quality-of-life", "physical problem", "cognitive", "affective/behavioral"]}]'
$ws.Range("G17").Value = '[{"factor": ["csf", "cerebrospinal fluid", "amyloid-beta1-42 (abeta42)", "concentration"], "outcome": ["mortality"]}, {"factor": ["amyloid-beta1-42 (abeta42)", "concentration", "plasma"], "outcome": ["mortality"]}, {"factor": ["change", "csf", "concentration", "abeta42"], "outcome": ["neurological status"]}]'
$ws.Range("G18").Value = '[{"factor": ["plasminogen activator receptor", "urokinase", "soluble", "supar"], "outcome": ["tbi", "traumatic brain injury"]}, {"factor": ["score", "glasgow coma scale"], "outcome": ["severity", "tbi"]}, {"factor": ["d-dimer"], "outcome": ["patient", "tbi", "mortality"]}]'
$ws.Range("G19").Value = '[{"factor": ["motor score", "gcs"], "outcome": ["mortality", "month"]}, {"factor": ["age"], "outcome": ["mortality", "month"]}, {"factor": ["sex"], "outcome": ["mortality", "month"]}, {"factor": ["mechanism", "injury"], "outcome": ["mortality", "month"]}, {"factor": ["glasgow coma scale"], "outcome": ["mortality", "month"]}, {"factor": ["intubation"], "outcome": ["mortality", "month"]}, {"factor": ["pupil"], "outcome": ["mortality", "month"]}, {"factor": ["systolic blood pressure"], "outcome": ["mortality", "month"]}, {"factor": ["respiratory rate"], "outcome": ["mortality", "month"]}, {"factor": ["body temperature"], "outcome": ["mortality", "month"]}, {"factor": ["ph", "arterial"], "outcome": ["mortality", "month"]}, {"factor": ["arterial partial pressure", "carbon dioxide"], "outcome": ["mortality", "month"]}, {"factor": ["arterial partial pressure"], "outcome": ["mortality", "month"]}, {"factor": ["serum sodium"], "outcome": ["mortality", "month"]}, {"factor": ["serum potassium"], "outcome": ["mortality", "month"]}, {"factor": ["serum chloride"], "outcome": ["mortality", "month"]}, {"factor": ["serum calcium"], "outcome": ["mortality", "month"]}, {"factor": ["serum glucose"], "outcome": ["mortality", "month"]}, {"factor": ["urea nitrogen"], "outcome": ["mortality", "month"]}, {"factor": ["creatinine"], "outcome": ["mortality", "month"]}, {"factor": ["international", "ratio"], "outcome": ["mortality", "month"]}]'
$ws.Range("G20").Value = '[{"factor": ["aptt"], "outcome": ["surgery", "deterioration"]}, {"factor": ["fdp"], "outcome": ["surgery", "deterioration"]}, {"factor": ["d-dimer"], "outcome": ["surgery", "deterioration"]}]'
$ws.Range("G21").Value = '[{"factor": ["central conduction time", "cct"], "outcome": ["long-term", "clinical outcome"]}, {"factor": ["latency"], "outcome": ["long-term", "clinical outcome"]}]'
$ws.Range("G24").Value = '[{"factor": ["tsp-1"], "outcome": ["unfavorable", "1-week", "mortality", "outcome"]}]'
$ws.Range("G25").Value = '[{"factor": ["plasma level", "brain-derived neurotrophic factor (bdnf)"], "outcome": ["tbi", "severe", "patient", "intensive care unit", "icu", "mortality"]}]'
$ws.Range("G26").Value = '[{"factor": ["crash-ct model"], "outcome": ["days", "death"]}, {"factor": ["age"], "outcome": ["older", "patient", "model", "performance"]}, {"factor": ["glasgow coma scale score"], "outcome": ["discrimination", "model"]}, {"factor": ["hosmer-lemeshow p value"], "outcome": ["model", "calibration"]}]'
$ws.Range("G27").Value = '[{"factor": ["time to death"], "outcome": ["withdrawal", "life-sustaining", "therapy"]}, {"factor": ["score", "glasgow coma scale"], "outcome": ["mortality"]}, {"factor": ["ais", "score", "head abbreviate injury scale"], "outcome": ["mortality"]}, {"factor": ["multiple", "comorbiditie"], "outcome": ["mortality"]}, {"factor": ["traumatic", "subarachnoid hemorrhage"], "outcome": ["mortality"]}, {"factor": ["intracerebral mass lesion"], "outcome": ["mortality"]}, {"factor": ["brainstem lesion"], "outcome": ["mortality"]}, {"factor": ["absent", "sign of compress", "basal cistern"], "outcome": ["mortality"]}]'
$ws.Range("G28").Value = '[{"factor": ["core model"], "outcome": ["mortality"]}, {"factor": ["extend model"], "outcome": ["mortality"]}, {"factor": ["lab model"], "outcome": ["mortality"]}, {"factor": ["core model"], "outcome": ["unfavorable", "outcome"]}, {"factor": ["extend model"], "outcome": ["unfavorable", "outcome"]}, {"factor": ["lab model"], "outcome": ["unfavorable", "outcome"]}]'
$ws.Range("G29").Value = '[{"factor": ["score", "glasgow coma scale"], "outcome": ["unfavorable", "1-week", "mortality", "outcome"]}]'
$ws.Range("G30").Value = '[{"factor": ["gcs", "glasgow coma scale"], "outcome": ["severity", "tbi"]}, {"factor": ["duration", "pta", "post-traumatic amnesia"], "outcome": ["tbi", "olfactory problem"]}]'
$ws.Range("G31").Value = '[{"factor": ["level", "il-6"], "outcome": ["development", "septic"]}, {"factor": ["c-reactive protein level"], "outcome": ["development", "multiple organ dysfunction"]}]'
$ws.Range("G32").Value = '[{"factor": ["rotterdam"], "outcome": ["hospital discharge", "death", "weeks"]}, {"factor": ["age"], "outcome": ["hospital discharge", "death", "weeks"]}, {"factor": ["sex"], "outcome": ["hospital discharge", "death", "weeks"]}, {"factor": ["glasgow coma scale score"], "outcome": ["hospital discharge", "death", "weeks"]}]'
$ws.Range("G33").Value = '[{"factor": ["admission", "serum", "glucose", "level of"], "outcome": ["patient", "outcome", "traumatic brain injury", "severe"]}]'
$ws.Range("G34").Value = '[{"factor": ["csf"], "outcome": ["gos score", "month"]}, {"factor": ["sfas"], "outcome": ["gos score", "month"]}, {"factor": ["il-10"], "outcome": ["gos score", "month"]}, {"factor": ["il-6"], "outcome": ["gos score", "month"]}, {"factor": ["svcam-1"], "outcome": ["gos score", "month"]}, {"factor": ["il-5"], "outcome": ["gos score", "month"]}, {"factor": ["il-8"], "outcome": ["gos score", "month"]}, {"factor": ["pc1"], "outcome": ["gos score", "month"]}, {"factor": [], "outcome": ["gos score", "month"]}, {"factor": [], "outcome": ["gos score", "month"]}]'
$ws.Range("G35").Value = '[{"factor": ["age"], "outcome": ["poor outcome"]}, {"factor": ["glasgow coma scale"], "outcome": ["poor outcome"]}, {"factor": ["severity score", "injury"], "outcome": ["poor outcome"]}, {"factor": ["ais", "head"], "outcome": ["poor outcome"]}]'
$ws.Range("G36").Value = '[{"factor": ["il-6"], "outcome": ["year", "gos", "favorable"]}, {"factor": ["pg"], "outcome": ["year", "gos", "favorable"]}, {"factor": ["gfap"], "outcome": ["year", "unfavorable", "score", "gos"]}, {"factor": ["pg"], "outcome": ["survival status", "year"]}, {"factor": ["gfap"], "outcome": ["survival status", "year"]}, {"factor": ["il-6"], "outcome": ["survival status", "year"]}]'
$ws.Range("G37").Value = '[{"factor": ["acute"], "outcome": ["score", "drs"]}, {"factor": ["subacute", "fa"], "outcome": ["score", "drs"]}]'
$ws.Range("G38").Value = '[{"factor": ["glasgow coma scale"], "outcome": ["mortality"]}, {"factor": ["mechanical ventilation"], "outcome": ["neurological"]}, {"factor": ["blood transfusion"], "outcome": ["neurological"]}, {"factor": ["neurosurgical intervention"], "outcome": ["neurological"]}, {"factor": ["injury", "concomitant"], "outcome": ["non-neurological", "complication"]}]'
